$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chargingdata")

# Rows 5:7 contain duplicate data (an exact repeat of rows 2:4) that needs to
# be removed. Deleting the entire rows shifts everything below up by three
# rows (old row 8 becomes the new row 5, old row 91 becomes the new row 88).
$ws.Rows("5:7").Delete() | Out-Null

# Set the final selection on the sheet, matching the author's last interaction.
$ws.Activate()
$ws.Range("I20").Select() | Out-Null
